# Fruta / hortaliza, semanal
# Insert a new weekly record at row 26 (pushing existing rows 26-34 down to
# 27-35), then populate the new row 26 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts rows 26..34 down to 27..35
# and carries the existing row formatting (e.g. the date style on column D).
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 44813
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112026
$ws.Range("G26").Value = "Haba"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 11000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 11500
$ws.Range("N26").Value = "`$/saco 25 kilos"
$ws.Range("O26").Value = "Región de Coquimbo"
$ws.Range("P26").Value = 460
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"
